$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'25.999.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = "'1.639.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.70%  '
$ws.Range("D5").Value = "'214.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.51%  '
$ws.Range("D8").Value = "'0.2584"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.45%  '
$ws.Range("D9").Value = "'0.06355"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.81%  '
$ws.Range("D10").Value = "'19.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.37%  '
$ws.Range("D11").Value = "'0.07738"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.63%  '
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("D13").Value = "'1.636.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("D14").Value = "'0.5478"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("D15").Value = "'0.0₅7744"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.56%  '
$ws.Range("D16").Value = "'64.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.58%  '
$ws.Range("D17").Value = "'26.020.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").Value = "'1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("D19").Value = "'196.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.89%  '
$ws.Range("D20").Value = "'4.461"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("D21").Value = "'9.966"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("D22").Value = "'6.128"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.42%  '
$ws.Range("D23").Value = "'1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.59%  '
$ws.Range("D24").Value = "'1.888"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = "'142.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.71%  '
$ws.Range("E26").Value = '  +9.93%  '
$ws.Range("D27").Value = "'6.872"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("D28").Value = "'15.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.87%  '
$ws.Range("D29").Value = "'1.241"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").Value = "'0.04881"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.65%  '
$ws.Range("D31").Value = "'3.284"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("D32").Value = "'3.215"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.65%  '
$ws.Range("D33").Value = "'1.555"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.87%  '
$ws.Range("D34").Value = "'2.374"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("D35").Value = "'0.9186"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.61%  '
$ws.Range("D36").Value = "'2.570"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.61%  '
$ws.Range("D37").Value = "'0.5548"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("D38").Value = "'1.106.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.28%  '
$ws.Range("D39").Value = "'0.01571"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.85%  '
$ws.Range("D40").Value = "'1.000"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.73%  '
$ws.Range("D41").Value = "'5.618"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.38%  '
$ws.Range("D42").Value = "'0.8049"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.48%  '
$ws.Range("D43").Value = "'98.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.94%  '
$ws.Range("E44").Value = '  -4.92%  '
$ws.Range("D45").Value = "'1.781.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("D47").Value = "'55.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.71%  '
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("D49").Value = "'0.05189"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.93%  '
$ws.Range("D50").Value = "'7.585"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.27%  '
$ws.Range("D51").Value = "'1.004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.21%  '
